# Version 3.0 Build 3031
# Update default configuration of
#   flexible_days_before_admission_for_CO  (row 5, column B)
#   flexibile_days_after_discharge_for_HO  (row 6, column B)
# from [0,0] to [-2,1]

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("setting_parameters")

$ws.Range("B5").Value = -2
$ws.Range("B6").Value = 1
